$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.294.12'
$ws.Range('E2').Value = '  +0.66%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.862.28'
$ws.Range('E3').Value = '  +0.74%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.18%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7002'
$ws.Range('E5').Value = '  -0.93%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '237.82'
$ws.Range('E6').Value = '  -0.20%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.08197'
$ws.Range('E8').Value = '  +9.49%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3045'
$ws.Range('E9').Value = '  -0.31%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.32'
$ws.Range('E10').Value = '  -0.36%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08172'
$ws.Range('E11').Value = '  +0.49%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.856.04'
$ws.Range('E12').Value = '  +0.47%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7175'
$ws.Range('E13').Value = '  -1.11%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.182'
$ws.Range('E14').Value = '  -0.90%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.40'
$ws.Range('E15').Value = '  +0.18%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.304.26'
$ws.Range('E16').Value = '  +0.56%  '

$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.786'
$ws.Range('E17').Value = '  -0.05%  '

$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007871'
$ws.Range('E18').Value = '  +2.57%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.40'
$ws.Range('E19').Value = '  +2.51%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '237.86'
$ws.Range('E20').Value = '  -0.89%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9998'
$ws.Range('E21').Value = '  +0.01%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.107.36'
$ws.Range('E22').Value = '  +0.67%  '

$ws.Range('E23').Value = '  +0.22%  '

$ws.Range('E24').Value = '  -1.17%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.31'
$ws.Range('E25').Value = '  +0.87%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.994'
$ws.Range('E26').Value = '  +0.22%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1441'
$ws.Range('E27').Value = '  -1.54%  '

$ws.Range('E28').Value = '  +0.64%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.983'
$ws.Range('E29').Value = '  +2.20%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.431'
$ws.Range('E30').Value = '  +3.34%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.433'
$ws.Range('E31').Value = '  -3.23%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.486'
$ws.Range('E32').Value = '  -0.53%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.064'
$ws.Range('E33').Value = '  +1.38%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05219'
$ws.Range('E34').Value = '  +0.98%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.173'
$ws.Range('E35').Value = '  -1.21%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7058'
$ws.Range('E36').Value = '  -0.23%  '

$ws.Range('E37').Value = '  -2.73%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.664'
$ws.Range('E38').Value = '  +0.86%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01853'
$ws.Range('E39').Value = '  -0.66%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.730'
$ws.Range('E40').Value = '  +2.04%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.144.92'
$ws.Range('E41').Value = '  +7.16%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9217'
$ws.Range('E42').Value = '  +1.13%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.978'
$ws.Range('E43').Value = '  -0.31%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4283'
$ws.Range('E44').Value = '  -0.34%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '70.85'
$ws.Range('E45').Value = '  +1.04%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('E46').Value = '  +0.10%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.99'
$ws.Range('E47').Value = '  +0.67%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.775'
$ws.Range('E48').Value = '  +1.14%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.004.74'
$ws.Range('E49').Value = '  +0.97%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.192'
$ws.Range('E50').Value = '  -0.14%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.976'
$ws.Range('E51').Value = '  -1.34%  '
